# "aggiornamento fino a 02/05" - append 6 more daily rows (27 Apr - 2 May 2021)
# after the existing data, continuing the same date/format pattern as the
# last existing row (row 238, date serial 44312).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 239
$startSerial = 44313
$numNewRows = 6

for ($i = 0; $i -lt $numNewRows; $i++) {
    $row = $startRow + $i
    $serial = $startSerial + $i

    # Carry the date cell's format (style "2": date number format, border,
    # bold font, centered alignment) down from the previous row, same as
    # dragging the fill handle in Excel would do.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 1
    $ws.Cells.Item($row, 4).Value = 46.70714619336758
}
